$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.152.95"
$ws.Range("E2").Value = "  +5.36%  "
$ws.Range("D3").Value = "2.238.69"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.37%  "
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.09"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.64%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.615"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0934"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.97%  "
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").Value = "2.576.72"
$ws.Range("E15").Value = "  +3.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.46%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.236.79"
$ws.Range("E17").Value = "  +3.10%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.813"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "43.059.49"
$ws.Range("E19").Value = "  +5.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000104"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "230.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.07%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("E28").Value = "  -5.06%  "
$ws.Range("E29").Value = "  +2.57%  "
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +22.79%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "174.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0796"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.122"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("E37").Value = "  +7.65%  "
$ws.Range("E38").Value = "  +5.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0334"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +16.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "13.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.200"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.45%  "
$ws.Range("E47").Value = "  +2.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.450"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +23.66%  "
$ws.Range("E49").Value = "  +2.59%  "
$ws.Range("E50").Value = "  +4.56%  "
$ws.Range("E51").Value = "  +2.80%  "
